$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.281.53"
$ws.Range("E2").Value = "  -0.14%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.186.15"
$ws.Range("E3").Value = "  -1.54%  "

# Row 5 - BNB
$ws.Range("D5").Value = "256.15"
$ws.Range("E5").Value = "  +5.00%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7 - Solana
$ws.Range("D7").Value = "68.15"
$ws.Range("E7").Value = "  -2.12%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.09%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.572"
$ws.Range("E9").Value = "  +2.71%  "

# Row 10 - was OKB, now Avalanche
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "36.98"
$ws.Range("E10").Value = "  -6.30%  "

# Row 11 - was Avalanche, now OKB
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "58.80"
$ws.Range("E11").Value = "  +1.21%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "0.0935"
$ws.Range("E12").Value = "  -1.96%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "7.00"
$ws.Range("E13").Value = "  +3.90%  "

# Row 14 - TRON
$ws.Range("D14").Value = "0.103"
$ws.Range("E14").Value = "  -0.82%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.508.04"
$ws.Range("E15").Value = "  -1.69%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.866"
$ws.Range("E16").Value = "  +2.79%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "14.36"
$ws.Range("E17").Value = "  -2.80%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.180.60"
$ws.Range("E18").Value = "  -1.73%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "41.197.35"
$ws.Range("E19").Value = "  -0.24%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +0.52%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.16"
$ws.Range("E21").Value = "  +1.13%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "71.92"
$ws.Range("E22").Value = "  -0.41%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "232.56"
$ws.Range("E23").Value = "  +0.28%  "

# Row 24 - ImmutableX
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  -3.10%  "

# Row 25 - was WEMIXToken, now Cosmos
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "11.70"
$ws.Range("E25").Value = "  +19.32%  "

# Row 26 - was Cosmos, now WEMIXToken
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "3.85"
$ws.Range("E26").Value = "  +6.29%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.13%  "

# Row 28 - PancakeSwap
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  +4.11%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  -5.53%  "

# Row 30 - Monero
$ws.Range("D30").Value = "168.92"
$ws.Range("E30").Value = "  -1.91%  "

# Row 31 - Toncoin
$ws.Range("E31").Value = "  -6.46%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "20.62"
$ws.Range("E32").Value = "  +0.66%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  -1.83%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.0747"
$ws.Range("E34").Value = "  +4.32%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  -0.63%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("D36").Value = "5.44"
$ws.Range("E36").Value = "  +3.94%  "

# Row 37 - InjectiveProtocol
$ws.Range("D37").Value = "26.25"
$ws.Range("E37").Value = "  +7.98%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +6.41%  "

# Row 39 - Filecoin
$ws.Range("E39").Value = "  -0.56%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +7.14%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  -3.81%  "

# Row 42 - Celestia
$ws.Range("D42").Value = "12.20"
$ws.Range("E42").Value = "  +12.96%  "

# Row 43 - THORChain
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").Value = "  -3.45%  "

# Row 44 - MultiversX
$ws.Range("D44").Value = "63.42"
$ws.Range("E44").Value = "  -3.58%  "

# Row 45 - FTXToken
$ws.Range("E45").Value = "  -0.57%  "

# Row 46 - Algorand
$ws.Range("E46").Value = "  -2.95%  "

# Row 47 - FraxShare
$ws.Range("D47").Value = "8.59"
$ws.Range("E47").Value = "  -2.53%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +0.14%  "

# Row 49 - was BinanceUSD, now ARBITRUM
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").Value = "  +8.25%  "

# Row 50 - was ARBITRUM, now BinanceUSD
$ws.Range("B50").Value = "BinanceUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.06%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  -0.31%  "
